# Applies the "first code review" changes to the Settings sheet of Config.xlsx:
#  - Row 11 (OutlookInputEmailAddress + mailto hyperlink) is replaced by a new
#    "OutputReportSheetName" row (hyperlink removed).
#  - A new "ColumnOfInterest" row is inserted at row 15.
#  - The previously-blank row 17 becomes the (now plain, non-hyperlinked)
#    "OutlookInputEmailAddress" row.
#  - Rows 19-32 are populated with the new Output-Report / Outlook / API
#    settings that were previously blank placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Remove the existing mailto: hyperlink on B11 (OutlookInputEmailAddress) ---
$ws.Hyperlinks.Delete()

# --- Row 11: OutlookInputEmailAddress -> OutputReportSheetName -------------
$ws.Range("A11").Value = "OutputReportSheetName"
$ws.Range("B11").Value = "Sheet1"
$ws.Range("C11").ClearContents()

# --- Row 15: new ColumnOfInterest row ---------------------------------------
$ws.Range("A15").Value = "ColumnOfInterest"
$ws.Range("B15").Value = "Animal "
$ws.Range("C15").Value = "It will search in queue item for this specific content."

# --- Row 17: OutlookInputEmailAddress (now without hyperlink) --------------
$ws.Range("A17").Value = "OutlookInputEmailAddress"
$ws.Range("B17").Value = "darius.dangi@fwfcompany.com"
$ws.Range("C17").Value = "Outlook email that is used to send mail messages."

# --- Row 19: GetRowIndexColumnName ------------------------------------------
$ws.Range("A19").Value = "GetRowIndexColumnName"
$ws.Range("B19").Value = "Name"
$ws.Range("C19").Value = "Searches for this specific column name to return the index of row."

# --- Row 20: OutputReportUpdatedColumn --------------------------------------
$ws.Range("A20").Value = "OutputReportUpdatedColumn"
$ws.Range("B20").Value = "Status"
$ws.Range("C20").Value = "This column value will be updated in the given dataTable."

# --- Row 22: OutlookInputAccount --------------------------------------------
$ws.Range("A22").Value = "OutlookInputAccount"
$ws.Range("B22").Value = "darius.dangi@fwfcompany.com"
$ws.Range("C22").Value = "This email address will be used to send the output report."

# --- Row 24: OutputReportSendTo ---------------------------------------------
$ws.Range("A24").Value = "OutputReportSendTo"
$ws.Range("B24").Value = "darius.dangi@fwfcompany.com"

# --- Row 25: OutputReportSubject --------------------------------------------
$ws.Range("A25").Value = "OutputReportSubject"
$ws.Range("B25").Value = "Output Report"

# --- Row 26: OutputReportBody -----------------------------------------------
$ws.Range("A26").Value = "OutputReportBody"
$ws.Range("B26").Value = "You have attached below the output report."

# --- Row 28: EndpointDogAPI --------------------------------------------------
$ws.Range("A28").Value = "EndpointDogAPI"
$ws.Range("B28").Value = "https://dog.ceo/api/breeds/image/random"
$ws.Range("C28").Value = "It is used to get a random picture of dogs."

# --- Row 29: EndpointCatAPI --------------------------------------------------
$ws.Range("A29").Value = "EndpointCatAPI"
$ws.Range("B29").Value = "https://api.thecatapi.com/v1/images/search"
$ws.Range("C29").Value = "It is used to get a random picture of cats."

# --- Row 31: GetAnimalFail ----------------------------------------------------
$ws.Range("A31").Value = "GetAnimalFail"
$ws.Range("B31").Value = "Failed to get success status from API call."

# --- Row 32: GetAnimalDownloadFail --------------------------------------------
$ws.Range("A32").Value = "GetAnimalDownloadFail"
$ws.Range("B32").Value = "Failed to download animal picture from API call."

# --- Move the active selection to A2 (matches the reviewed workbook state) ---
$ws.Activate()
$ws.Range("A2").Select()
